$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCT1")
Write-Host $ws.Name
